$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("F37").Value = 'Marines'
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 'AS Kigali'
$ws.Range("J37").Value = 2.99
$ws.Range("L37").Value = 2.75
$ws.Range("M37").Value = '11/10/2023 14:57'
$ws.Range("N37").Value = 2.78
$ws.Range("P37").Value = 2.88
$ws.Range("Q37").Value = '11/10/2023 14:57'
$ws.Range("R37").Value = 2.21
$ws.Range("T37").Value = 2.54
$ws.Range("U37").Value = '11/10/2023 14:57'
$ws.Range("V37").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/marines-as-kigali/QHn271so/'

# Row 38
$ws.Range("F38").Value = 'Police'
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 'Muhazi United'
$ws.Range("J38").Value = 1.53
$ws.Range("L38").Value = 1.69
$ws.Range("M38").Value = '11/10/2023 14:56'
$ws.Range("N38").Value = 3.41
$ws.Range("P38").Value = 3.35
$ws.Range("Q38").Value = '11/10/2023 14:56'
$ws.Range("R38").Value = 4.88
$ws.Range("T38").Value = 4.66
$ws.Range("U38").Value = '11/10/2023 14:56'
$ws.Range("V38").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/police-muhazi-united/GrcSD34H/'

# Row 47
$ws.Range("F47").Value = 'AS Kigali'
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 'Police'
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 1.97
$ws.Range("L47").Value = 2.07
$ws.Range("M47").Value = '15/10/2023 11:02'
$ws.Range("N47").Value = 2.82
$ws.Range("P47").Value = 2.81
$ws.Range("Q47").Value = '15/10/2023 13:01'
$ws.Range("R47").Value = 3.49
$ws.Range("T47").Value = 3.7
$ws.Range("U47").Value = '15/10/2023 11:02'
$ws.Range("V47").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/as-kigali-police/KzfJGKxk/'

# Row 48
$ws.Range("F48").Value = 'Musanze'
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 'Rayon Sport'
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 3.66
$ws.Range("L48").Value = 3.04
$ws.Range("M48").Value = '15/10/2023 14:58'
$ws.Range("N48").Value = 2.93
$ws.Range("P48").Value = 2.4
$ws.Range("Q48").Value = '15/10/2023 14:58'
$ws.Range("R48").Value = 1.87
$ws.Range("T48").Value = 2.79
$ws.Range("U48").Value = '15/10/2023 14:58'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-rayon-sport/bRhNFvhe/'

# Row 59
$ws.Range("F59").Value = 'Sunrise'
$ws.Range("H59").Value = 'Muhazi United'
$ws.Range("I59").Value = 2
$ws.Range("J59").Value = 2.01
$ws.Range("L59").Value = 2.05
$ws.Range("M59").Value = '29/10/2023 04:30'
$ws.Range("N59").Value = 2.82
$ws.Range("P59").Value = 2.92
$ws.Range("Q59").Value = '29/10/2023 12:02'
$ws.Range("R59").Value = 3.39
$ws.Range("T59").Value = 3.59
$ws.Range("U59").Value = '29/10/2023 04:30'
$ws.Range("V59").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/sunrise-muhazi-united/hO8zty60/'

# Row 61
$ws.Range("F61").Value = 'Musanze'
$ws.Range("H61").Value = 'AS Kigali'
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2.78
$ws.Range("L61").Value = 2.97
$ws.Range("M61").Value = '29/10/2023 13:46'
$ws.Range("N61").Value = 2.6
$ws.Range("P61").Value = 2.53
$ws.Range("Q61").Value = '29/10/2023 13:46'
$ws.Range("R61").Value = 2.54
$ws.Range("T61").Value = 2.69
$ws.Range("U61").Value = '29/10/2023 13:46'
$ws.Range("V61").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/musanze-as-kigali/6RCvueL6/'

# Row 68
$ws.Range("F68").Value = 'Amagaju'
$ws.Range("H68").Value = 'Police'
$ws.Range("I68").Value = 2
$ws.Range("J68").Value = 2.53
$ws.Range("K68").Value = '25/11/2023 12:44'
$ws.Range("L68").Value = 2.93
$ws.Range("M68").Value = '25/11/2023 13:16'
$ws.Range("N68").Value = 2.63
$ws.Range("O68").Value = '25/11/2023 12:44'
$ws.Range("P68").Value = 2.63
$ws.Range("Q68").Value = '25/11/2023 12:44'
$ws.Range("R68").Value = 3.04
$ws.Range("S68").Value = '25/11/2023 12:44'
$ws.Range("T68").Value = 2.61
$ws.Range("U68").Value = '25/11/2023 13:16'
$ws.Range("V68").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/amagaju-police/z7QQjhKI/'

# Row 69
$ws.Range("F69").Value = 'Etincelles'
$ws.Range("H69").Value = 'Rayon Sport'
$ws.Range("I69").Value = 1
$ws.Range("J69").Value = 4.18
$ws.Range("K69").Value = '11/11/2023 03:13'
$ws.Range("L69").Value = 3.12
$ws.Range("M69").Value = '25/11/2023 13:56'
$ws.Range("N69").Value = 3.15
$ws.Range("O69").Value = '11/11/2023 03:13'
$ws.Range("P69").Value = 2.74
$ws.Range("Q69").Value = '25/11/2023 13:56'
$ws.Range("R69").Value = 1.77
$ws.Range("S69").Value = '11/11/2023 03:13'
$ws.Range("T69").Value = 2.38
$ws.Range("U69").Value = '25/11/2023 13:56'
$ws.Range("V69").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/etincelles-rayon-sport/tpRMiY4C/'

# Row 70
$ws.Range("F70").Value = 'Bugesera'
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 'Marines'
$ws.Range("J70").Value = 1.99
$ws.Range("L70").Value = 2.4
$ws.Range("M70").Value = '25/11/2023 13:54'
$ws.Range("N70").Value = 3.07
$ws.Range("P70").Value = 2.91
$ws.Range("Q70").Value = '25/11/2023 13:54'
$ws.Range("R70").Value = 3.42
$ws.Range("T70").Value = 2.9
$ws.Range("U70").Value = '25/11/2023 13:54'
$ws.Range("V70").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/bugesera-marines/2gSIhEk6/'

# Row 71
$ws.Range("F71").Value = 'APR'
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 'AS Kigali'
$ws.Range("J71").Value = 1.7
$ws.Range("L71").Value = 1.77
$ws.Range("M71").Value = '25/11/2023 12:54'
$ws.Range("N71").Value = 3.13
$ws.Range("P71").Value = 3.04
$ws.Range("Q71").Value = '25/11/2023 12:54'
$ws.Range("R71").Value = 4.63
$ws.Range("T71").Value = 4.74
$ws.Range("U71").Value = '25/11/2023 12:54'
$ws.Range("V71").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/apr-as-kigali/nqgfZPlo/'

# Row 72
$ws.Range("F72").Value = 'Kiyovu'
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 'Sunrise'
$ws.Range("J72").Value = 2.33
$ws.Range("K72").Value = '25/11/2023 13:20'
$ws.Range("L72").Value = 2.33
$ws.Range("M72").Value = '25/11/2023 13:20'
$ws.Range("O72").Value = '25/11/2023 13:20'
$ws.Range("P72").Value = 2.73
$ws.Range("Q72").Value = '25/11/2023 13:20'
$ws.Range("R72").Value = 3.23
$ws.Range("S72").Value = '25/11/2023 13:20'
$ws.Range("T72").Value = 3.23
$ws.Range("U72").Value = '25/11/2023 13:20'
$ws.Range("V72").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/kiyovu-sunrise/llXxb9Ro/'

# Row 73
$ws.Range("F73").Value = 'Gasogi United'
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 'Musanze'
$ws.Range("J73").Value = 2.62
$ws.Range("K73").Value = '12/11/2023 03:12'
$ws.Range("L73").Value = 2.78
$ws.Range("M73").Value = '26/11/2023 10:04'
$ws.Range("N73").Value = 2.73
$ws.Range("O73").Value = '12/11/2023 03:12'
$ws.Range("P73").Value = 2.75
$ws.Range("Q73").Value = '26/11/2023 12:05'
$ws.Range("R73").Value = 2.68
$ws.Range("S73").Value = '12/11/2023 03:12'
$ws.Range("T73").Value = 2.58
$ws.Range("U73").Value = '26/11/2023 10:04'
$ws.Range("V73").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/gasogi-united-musanze/SvWtcTth/'

# Row 74
$ws.Range("E74").Value = 45256.58333333334
$ws.Range("F74").Value = 'Mukura Victory Sports'
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 'Muhazi United'
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1.8
$ws.Range("K74").Value = '12/11/2023 13:03'
$ws.Range("L74").Value = 1.8
$ws.Range("M74").Value = '12/11/2023 13:03'
$ws.Range("N74").Value = 3.02
$ws.Range("O74").Value = '12/11/2023 13:03'
$ws.Range("P74").Value = 3.02
$ws.Range("Q74").Value = '12/11/2023 13:03'
$ws.Range("R74").Value = 4.57
$ws.Range("S74").Value = '12/11/2023 13:03'
$ws.Range("T74").Value = 4.57
$ws.Range("U74").Value = '12/11/2023 13:03'
$ws.Range("V74").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/mukura-victory-sports-muhazi-united/z5Vpdmdb/'

# New rows 75 and 76 - copy formatting from row 74, then populate values
$ws.Range("A74:V74").Copy()
$ws.Range("A75:V76").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 75
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = 'rwanda'
$ws.Range("C75").Value = 'premier-league'
$ws.Range("D75").Value = '2023-2024'
$ws.Range("E75").Value = 45258.58333333334
$ws.Range("F75").Value = 'Police'
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 'Rayon Sport'
$ws.Range("I75").Value = 2
$ws.Range("J75").Value = 2.73
$ws.Range("K75").Value = '28/11/2023 03:12'
$ws.Range("L75").Value = 2.6
$ws.Range("M75").Value = '28/11/2023 13:45'
$ws.Range("N75").Value = 2.76
$ws.Range("O75").Value = '28/11/2023 03:12'
$ws.Range("P75").Value = 2.61
$ws.Range("Q75").Value = '28/11/2023 13:45'
$ws.Range("R75").Value = 2.62
$ws.Range("S75").Value = '28/11/2023 03:12'
$ws.Range("T75").Value = 2.97
$ws.Range("U75").Value = '28/11/2023 13:45'
$ws.Range("V75").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/police-rayon-sport/KMAts8XA/'

# Row 76
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = 'rwanda'
$ws.Range("C76").Value = 'premier-league'
$ws.Range("D76").Value = '2023-2024'
$ws.Range("E76").Value = 45259.58333333334
$ws.Range("F76").Value = 'Sunrise'
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 'APR'
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 5.7
$ws.Range("K76").Value = '29/11/2023 13:33'
$ws.Range("L76").Value = 6.34
$ws.Range("M76").Value = '29/11/2023 13:56'
$ws.Range("N76").Value = 3.65
$ws.Range("O76").Value = '29/11/2023 13:33'
$ws.Range("P76").Value = 3.45
$ws.Range("Q76").Value = '29/11/2023 13:56'
$ws.Range("R76").Value = 1.45
$ws.Range("S76").Value = '29/11/2023 13:33'
$ws.Range("T76").Value = 1.52
$ws.Range("U76").Value = '29/11/2023 13:56'
$ws.Range("V76").Value = 'https://www.betexplorer.com/football/rwanda/premier-league/sunrise-apr/6J6xrlI4/'
